$wb = $excel.ActiveWorkbook

# TC004 is the sheet that should hold the "Yellow" value and become the active tab
$ws4 = $wb.Worksheets.Item("TC004")

# Update the Color cell for the "Dresses" row from "Blue" to "Yellow"
$ws4.Range("E2").Value = "Yellow"

# Select the changed cell and activate this sheet so it becomes the active tab
$ws4.Activate()
$ws4.Range("E2").Select()
